# Day 2 slides were added to the deck; as part of that commit the two
# theme parts (ppt/theme/theme1.xml and ppt/theme/theme2.xml) swapped
# places: the "Simple Light" theme moved into theme1.xml and the plain
# "Default" colour theme moved into theme2.xml (the theme that is
# actually wired to the slide master, i.e. the one driving every
# slide's colours).
#
# The PowerPoint object model doesn't let us rename/relink OOXML parts,
# but it does let us rewrite the 12 theme colour slots on the live
# theme (Master.Theme.ThemeColorScheme), which is the part of the swap
# that actually affects how the deck looks. Drive the currently active
# theme's colours to the "Default" scheme's RGB values.

function HexToRgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# "Default" clrScheme, in MsoThemeColorSchemeIndex order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$defaultScheme = @(
    "000000",
    "FFFFFF",
    "158158",
    "F3F3F3",
    "058DC7",
    "50B432",
    "ED561B",
    "EDEF00",
    "24CBE5",
    "64E572",
    "2200CC",
    "551A8B"
)

$p = $ppt.ActivePresentation
$theme = $p.SlideMaster.Theme
$colorScheme = $theme.ThemeColorScheme

for ($i = 1; $i -le $defaultScheme.Count; $i++) {
    $colorScheme.Colors($i).RGB = HexToRgb $defaultScheme[$i - 1]
}
